$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fullname (column C) for existing rows 2-5 to "Phạm Thanh Hà"
$ws.Range("C2").Value = "Phạm Thanh Hà"
$ws.Range("C3").Value = "Phạm Thanh Hà"
$ws.Range("C4").Value = "Phạm Thanh Hà"
$ws.Range("C5").Value = "Phạm Thanh Hà"

# Update Content (column K) for row 3
$ws.Range("K3").Value = "Tình một đêm"

# Add new row 6 with a new certificate record
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "HE130576"
$ws.Range("C6").Value = "Phạm Thanh Hà"
$ws.Range("D6").Value = "26/09/1999"
$ws.Range("E6").Value = "Nam"
$ws.Range("F6").Value = "Kinh"
$ws.Range("G6").Value = "Việt Nam"
$ws.Range("H6").Value = "Hà Nội"
$ws.Range("I6").Value = "ABC101"
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = "Học"
$ws.Range("L6").Value = "ĐH200305"

# Match the saved selection/active cell from the authored edit
$ws.Range("H23").Select() | Out-Null
